$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header B1 from "domain" to "domain_external_id"
# (Business Concept file manager domain name field)
$ws.Range("B1").Value = "domain_external_id"

# Move the active selection to the edited header cell
[void]$ws.Range("B1").Select()
